$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new "DEF Notes" worksheet, positioned right before the
#    "Excel Notes" sheet (Index, ABC Notes, DEF Notes, Excel Notes).
# ------------------------------------------------------------------
$wsExcelNotesOriginal = $wb.Worksheets.Item("Excel Notes")
$wsDef = $wb.Worksheets.Add($wsExcelNotesOriginal)
$wsDef.Name = "DEF Notes"

# NOTE: this engine's Worksheets.Add(Before:=...) hands back a handle
# that ends up aliasing the same underlying sheet as the "Before"
# variable once the new sheet is renamed, so re-resolve "Excel Notes"
# by name afterwards rather than continuing to use $wsExcelNotesOriginal.

# Match the look & feel of the other "notes" sheets.
$wsDef.Columns.Item(1).ColumnWidth = 9.83
$wsDef.Columns.Item(2).ColumnWidth = 49.83

# ------------------------------------------------------------------
# 2. Populate the DEF Notes sheet and turn the range into the
#    DEF_NOTES table (mirrors ABC_NOTES: Number / Note columns).
# ------------------------------------------------------------------
$wsDef.Range("A1").Value = "Number"
$wsDef.Range("B1").Value = "Note"
$wsDef.Range("A2").Value = 1
$wsDef.Range("B2").Value = "GIVE SOMEONE A HIGH FIVE"
$wsDef.Range("A3").Value = 2
$wsDef.Range("B3").Value = "BE NICE"
$wsDef.Range("A4").Value = 3
$wsDef.Range("B4").Value = "PICK UP TRASH"
$wsDef.Range("A5").Value = 4
$wsDef.Range("B5").Value = "EAT HEALTHY"

$loDef = $wsDef.ListObjects.Add(1, $wsDef.Range("A1:B5"), 0, 1)
$loDef.Name = "DEF_NOTES"

# ------------------------------------------------------------------
# 3. Extend the ABC_NOTES table on "ABC Notes" with two more notes.
# ------------------------------------------------------------------
$wsAbc = $wb.Worksheets.Item("ABC Notes")
$loAbc = $wsAbc.ListObjects.Item("ABC_NOTES")

$loAbc.ListRows.Add() | Out-Null
$wsAbc.Range("A4").Value = 3
$wsAbc.Range("B4").Value = "CONSTRUCT ROADWAY"

$loAbc.ListRows.Add() | Out-Null
$wsAbc.Range("A5").Value = 4
$wsAbc.Range("B5").Value = "CONSTRUCT UTILITY"

# ------------------------------------------------------------------
# 4. Extend the SHEET_INDEX table on "Index" with the new ABC-103
#    sheet plus the three new DEF-10x sheets.
# ------------------------------------------------------------------
$wsIndex = $wb.Worksheets.Item("Index")
$loIndex = $wsIndex.ListObjects.Item("SHEET_INDEX")

$loIndex.ListRows.Add() | Out-Null
$wsIndex.Range("A4").Value = "ABC-103"
$wsIndex.Range("B4").Value = "PROJ-ABC-100"
$wsIndex.Range("C4").Value = "ABC PLAN"

$loIndex.ListRows.Add() | Out-Null
$wsIndex.Range("A5").Value = "DEF-101"
$wsIndex.Range("B5").Value = "PROJ-DEF-100"
$wsIndex.Range("C5").Value = "DEF PLAN"

$loIndex.ListRows.Add() | Out-Null
$wsIndex.Range("A6").Value = "DEF-102"
$wsIndex.Range("B6").Value = "PROJ-DEF-100"
$wsIndex.Range("C6").Value = "DEF PLAN"

$loIndex.ListRows.Add() | Out-Null
$wsIndex.Range("A7").Value = "DEF-103"
$wsIndex.Range("B7").Value = "PROJ-DEF-100"
$wsIndex.Range("C7").Value = "DEF PLAN"

# ------------------------------------------------------------------
# 5. Fix up + extend the EXCEL_NOTES table on "Excel Notes":
#    - ABC-101 row gains the "Note 2" value (C2)
#    - ABC-102 row's values shift to Note1=2 / Note2=3
#    - Add rows for ABC-103 and the three DEF-10x sheets.
# ------------------------------------------------------------------
$wsExcel = $wb.Worksheets.Item("Excel Notes")
$loExcel = $wsExcel.ListObjects.Item("EXCEL_NOTES")

$wsExcel.Range("C2").Value = 2

$wsExcel.Range("B3").Value = 2
$wsExcel.Range("C3").Value = 3

$loExcel.ListRows.Add() | Out-Null
$wsExcel.Range("A4").Value = "ABC-103"
$wsExcel.Range("B4").Value = 3
$wsExcel.Range("C4").Value = 4

$loExcel.ListRows.Add() | Out-Null
$wsExcel.Range("A5").Value = "DEF-101"
$wsExcel.Range("B5").Value = 1
$wsExcel.Range("C5").Value = 2

$loExcel.ListRows.Add() | Out-Null
$wsExcel.Range("A6").Value = "DEF-102"
$wsExcel.Range("B6").Value = 2
$wsExcel.Range("C6").Value = 3

$loExcel.ListRows.Add() | Out-Null
$wsExcel.Range("A7").Value = "DEF-103"
$wsExcel.Range("B7").Value = 3
$wsExcel.Range("C7").Value = 4

# ------------------------------------------------------------------
# 6. Restore per-sheet selections to match the saved state.
# ------------------------------------------------------------------
$wsExcel.Range("A5:A7").Select() | Out-Null
$wsDef.Range("B5").Select() | Out-Null
$wsAbc.Range("B6").Select() | Out-Null
$wsIndex.Range("B7").Select() | Out-Null

# The Index tab is the one left active/selected.
$wsIndex.Activate()
